$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2-14): 46073 -> 46074
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 46074
}

# Update A/B/G for rows 7,8,9,10,12,13,14 (row 11 unchanged apart from C above)
$rowData = @{
    7  = @("A 62884-2021", 44504, 0.8)
    8  = @("A 28266-2025", 45818.56381944445, 1.9)
    9  = @("A 25015-2023", 45085.6989699074, 1.8)
    10 = @("A 19922-2025", 45771.63034722222, 10.1)
    12 = @("A 3402-2026", 46042.39047453704, 5.5)
    13 = @("A 14271-2021", 44278, 6.7)
    14 = @("A 25634-2025", 45803.59570601852, 6)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 7).Value = $vals[2]
}
